$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.45843057308078
$ws.Range("D2").Value = 5.673609554272452
$ws.Range("E2").Value = 13.82782050721514
$ws.Range("F2").Value = 34.00327528687052
$ws.Range("G2").Value = 3.640837487772983
$ws.Range("L2").Value = 9.225315554510207
$ws.Range("O2").Value = 29.721790951329
$ws.Range("C3").Value = 13.36616548088737
$ws.Range("D3").Value = 5.686663600238081
$ws.Range("E3").Value = 13.77440445474068
$ws.Range("F3").Value = 33.45556674027195
$ws.Range("G3").Value = 3.645260097764412
$ws.Range("L3").Value = 9.233124427980497
$ws.Range("O3").Value = 29.35654874916341
$ws.Range("C4").Value = 13.31305664959182
$ws.Range("D4").Value = 5.695707863816803
$ws.Range("E4").Value = 13.74516280287109
$ws.Range("F4").Value = 33.12610478851779
$ws.Range("G4").Value = 3.648115624878392
$ws.Range("L4").Value = 9.239670170994421
$ws.Range("O4").Value = 29.13942427315405
$ws.Range("C5").Value = 13.29232169284377
$ws.Range("D5").Value = 5.699651528292306
$ws.Range("E5").Value = 13.7341468655752
$ws.Range("F5").Value = 32.9937423111425
$ws.Range("G5").Value = 3.649314633341602
$ws.Range("L5").Value = 9.242777148845271
$ws.Range("O5").Value = 29.05283770159671
$ws.Range("C6").Value = 13.28893392007829
$ws.Range("D6").Value = 5.700321935737565
$ws.Range("E6").Value = 13.73237220620938
$ws.Range("F6").Value = 32.97188325159983
$ws.Range("G6").Value = 3.649515867633657
$ws.Range("L6").Value = 9.243319581575239
$ws.Range("O6").Value = 29.03857723123516
$ws.Range("C7").Value = 13.31277331633441
$ws.Range("D7").Value = 5.695760005535521
$ws.Range("E7").Value = 13.74501058573994
$ws.Range("F7").Value = 33.12431179315842
$ws.Range("G7").Value = 3.648131651766386
$ws.Range("L7").Value = 9.239710294240801
$ws.Range("O7").Value = 29.13824874237357
$ws.Range("C8").Value = 13.42589282295306
$ws.Range("D8").Value = 5.677896594854479
$ws.Range("E8").Value = 13.80866860424339
$ws.Range("F8").Value = 33.81311469682641
$ws.Range("G8").Value = 3.642333429923385
$ws.Range("L8").Value = 9.227644183275331
$ws.Range("O8").Value = 29.59443518759499
$ws.Range("C9").Value = 13.67501977545938
$ws.Range("D9").Value = 5.651062814987629
$ws.Range("E9").Value = 13.96139696018246
$ws.Range("F9").Value = 35.21014331858665
$ws.Range("G9").Value = 3.632067476928846
$ws.Range("L9").Value = 9.217907576451886
$ws.Range("O9").Value = 30.54108561839029
$ws.Range("C10").Value = 13.87351330624757
$ws.Range("D10").Value = 5.636387773005543
$ws.Range("E10").Value = 14.09010922011087
$ws.Range("F10").Value = 36.2541036030977
$ws.Range("G10").Value = 3.625188936926795
$ws.Range("L10").Value = 9.219279807822971
$ws.Range("O10").Value = 31.26204015738254
$ws.Range("C11").Value = 13.96689110188249
$ws.Range("D11").Value = 5.630815051330783
$ws.Range("E11").Value = 14.15211736002805
$ws.Range("F11").Value = 36.73061710149481
$ws.Range("G11").Value = 3.622201851554446
$ws.Range("L11").Value = 9.221760514494699
$ws.Range("O11").Value = 31.59417816948039
$ws.Range("C12").Value = 14.00266896467729
$ws.Range("D12").Value = 5.628864128346564
$ws.Range("E12").Value = 14.17608267467964
$ws.Range("F12").Value = 36.91111181689759
$ws.Range("G12").Value = 3.621090983327049
$ws.Range("L12").Value = 9.222966993544306
$ws.Range("O12").Value = 31.72043532706199
$ws.Range("C13").Value = 13.99494539814853
$ws.Range("D13").Value = 5.629277196091112
$ws.Range("E13").Value = 14.17089999553255
$ws.Range("F13").Value = 36.87223969149542
$ws.Range("G13").Value = 3.621329329158109
$ws.Range("L13").Value = 9.222695277311534
$ws.Range("O13").Value = 31.69322391484194
$ws.Range("C14").Value = 13.96982634539394
$ws.Range("D14").Value = 5.63065135006889
$ws.Range("E14").Value = 14.15407937401398
$ws.Range("F14").Value = 36.74546654308774
$ws.Range("G14").Value = 3.62211005420924
$ws.Range("L14").Value = 9.221854419159637
$ws.Range("O14").Value = 31.60455638953842
$ws.Range("C15").Value = 13.95449382897091
$ws.Range("D15").Value = 5.631513831998693
$ws.Range("E15").Value = 14.14383891528421
$ws.Range("F15").Value = 36.66781545417791
$ws.Range("G15").Value = 3.622590907166764
$ws.Range("L15").Value = 9.221374153538996
$ws.Range("O15").Value = 31.55030453097936
$ws.Range("C16").Value = 13.86747070858824
$ws.Range("D16").Value = 5.636774222147579
$ws.Range("E16").Value = 14.08612531612187
$ws.Range("F16").Value = 36.22297902924894
$ws.Range("G16").Value = 3.62538699510964
$ws.Range("L16").Value = 9.219155059316932
$ws.Range("O16").Value = 31.24040815154056
$ws.Range("C17").Value = 13.81485638975998
$ws.Range("D17").Value = 5.64028433899758
$ws.Range("E17").Value = 14.05159630043361
$ws.Range("F17").Value = 35.95036740577846
$ws.Range("G17").Value = 3.627138571631217
$ws.Range("L17").Value = 9.218269343641715
$ws.Range("O17").Value = 31.05128235415341
$ws.Range("C18").Value = 13.7848857124775
$ws.Range("D18").Value = 5.642407044972813
$ws.Range("E18").Value = 14.03206227613798
$ws.Range("F18").Value = 35.79373059941138
$ws.Range("G18").Value = 3.628159406073175
$ws.Range("L18").Value = 9.217934644857726
$ws.Range("O18").Value = 30.94290182568165
$ws.Range("C19").Value = 13.77478901584925
$ws.Range("D19").Value = 5.643143558397613
$ws.Range("E19").Value = 14.02550478367895
$ws.Range("F19").Value = 35.74072956936628
$ws.Range("G19").Value = 3.628507344417283
$ws.Range("L19").Value = 9.217851327196829
$ws.Range("O19").Value = 30.90627823508239
$ws.Range("C20").Value = 13.82042726717822
$ws.Range("D20").Value = 5.639899934923479
$ws.Range("E20").Value = 14.05523830757428
$ws.Range("F20").Value = 35.97937193093664
$ws.Range("G20").Value = 3.626950730202498
$ws.Range("L20").Value = 9.218345543099947
$ws.Range("O20").Value = 31.07137460360341
$ws.Range("C21").Value = 13.97719329052193
$ws.Range("D21").Value = 5.630243397387478
$ws.Range("E21").Value = 14.15900696467694
$ws.Range("F21").Value = 36.78270297203278
$ws.Range("G21").Value = 3.621880187184462
$ws.Range("L21").Value = 9.222094150543068
$ws.Range("O21").Value = 31.63058798430457
$ws.Range("C22").Value = 14.0820713647512
$ws.Range("D22").Value = 5.624861408286137
$ws.Range("E22").Value = 14.22964135465721
$ws.Range("F22").Value = 37.3079340442691
$ws.Range("G22").Value = 3.618684413146165
$ws.Range("L22").Value = 9.226100839348996
$ws.Range("O22").Value = 31.99883269342864
$ws.Range("C23").Value = 14.02588293528819
$ws.Range("D23").Value = 5.627648625605822
$ws.Range("E23").Value = 14.19168929947763
$ws.Range("F23").Value = 37.02764794139443
$ws.Range("G23").Value = 3.620379297130374
$ws.Range("L23").Value = 9.223819945605511
$ws.Range("O23").Value = 31.80207850291596
$ws.Range("C24").Value = 13.81790780568107
$ws.Range("D24").Value = 5.640073397965302
$ws.Range("E24").Value = 14.05359076728036
$ws.Range("F24").Value = 35.96625868700853
$ws.Range("G24").Value = 3.627035610225263
$ws.Range("L24").Value = 9.218310549773793
$ws.Range("O24").Value = 31.06228979734317
$ws.Range("C25").Value = 13.60481832645837
$ws.Range("D25").Value = 5.657440090125072
$ws.Range("E25").Value = 13.9171410040944
$ws.Range("F25").Value = 34.82835901985344
$ws.Range("G25").Value = 3.63472743262612
$ws.Range("L25").Value = 9.219046490638757
$ws.Range("O25").Value = 30.28005906737459
